$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = "`r`n"
$ws.Range("C21").Value = "Ханты-Мансийский автономный округ-Югра, `r`nг. Нижневартовск, Кузоваткина ул., 14"
$ws.Range("C22").Value = "Ханты-Мансийский автономный округ-Югра, `r`nг. Нижневартовск, ул. Ленина, д. 17, корп. П"
$ws.Range("C25").Value = "Ханты-Мансийский автономный округ-Югра,`r`n г. Нижневартовск,ул. Кузоватнкина, д. 14"
$ws.Range("C27").Value = "Ханты-Мансийский автономный округ-Югра,`r`n г. Нижневартовск, ул. Ленина, д. 4"
$ws.Range("C28").Value = "Ханты-Мансийский автономный округ-Югра, г. Нягань, `r`nул. Сибирская, д. 10, корп. 1`r`n"
$ws.Range("B30").Value = "ООО «Тагульское»`r`n"
$ws.Range("C30").Value = "Красноярский край, `r`nг. Красноярск, `r`nул. 78 Добровольческой бригады, д. 15`r`n"
$ws.Range("C36").Value = "Ханты-Мансийский автономный округ - Югра, город Когалым,улица Мира, дом 23, корпус А`r`n"
$ws.Range("C64").Value = "Ямало-Ненецкий автономный округ, г. Надым, ул. Комсомольская, д.16, кв. 36`r`n"
$ws.Range("B149").Value = "Акционерное общество `r`n«Нефтяная компания Дулисьма»`r`n(АО «НК Дулисьма»)`r`n"
$ws.Range("C152").Value = "Иркутская область, г. Иркутск, пр-т Большой Литейный, д. 4`r`n"
$ws.Range("C177").Value = "Республика Татарстан, р-н Черемшанский, с. Черемшан, `r`nул. Советская, д. 32, пом. 316`r`n"
$ws.Range("C192").Value = "край Пермский, город Пермь,`r`nулица Петропавловская, дом 123, офис 4`r`n"
$ws.Range("C213").Value = "Республика Татарстан, город Казань, улица Муштари, дом 2А, пом/офис 100Н/41`r`n"
$ws.Range("C239").Value = "Самарская область, г. Самара,`r`nул. Клиническая, д. 154, литера ИИ2И1, комн.2  `r`n"
$ws.Range("C283").Value = "Республика Татарстан, `r`nг. Альметьевск, ул. Маяковского, д. 116"
$ws.Range("B296").Value = "Акционерное общество «Томскнефть» Восточной Нефтяной Компании `r`n(АО «Томскнефть» ВНК)`r`n"
$ws.Range("C311").Value = "область Московская,`r`nрайон Дмитровский, город Дмитров, улица Профессиональная, дом 135, корпус 3, пом. 125`r`n"
$ws.Range("C334").Value = "Ханты-Мансийский Автономный округ – Югра, город Нижневартовск, улица Мира, дом 24`r`n"
$ws.Range("B376").Value = "ООО «ВостокИнвестНефть»`r`n"
$ws.Range("C376").Value = "Ульяновская область, р.п. Новоспасское, ул. Гагарина, `r`nд. 25 `r`n"
$ws.Range("B377").Value = "ООО «НК «ГНТ»`r`n"
$ws.Range("C377").Value = "Саратовская область, г. Саратов, ул. Соборная, д. 21М`r`n"
$ws.Range("C381").Value = "Тюменская область, г. Тюмень, ул. Республики, `r`nд. 143А, оф. 1401`r`n"
$ws.Range("B384").Value = "`r`nООО « Преображенское»`r`n"
$ws.Range("C384").Value = "Оренбургская область, г. Оренбург, ул. Джангильдина,`r`nд. 3, пом. 12"
$ws.Range("B385").Value = "АО «Нефтегазрезерв»`r`n"
$ws.Range("C385").Value = "Саратовская область,  г. Саратов, ул. им. Пугачева Е.И.,`r`nд. 159, оф. 905А`r`n"
$ws.Range("B386").Value = "ООО «Дубровинское»`r`n"
$ws.Range("C386").Value = "Удмуртская республика, г. Ижевск, ул. им. Репина, д. 35/1, кв. 106`r`n"
$ws.Range("B387").Value = "ООО «Азинское»`r`n"
$ws.Range("C387").Value = "Удмуртская республика, г. Ижевск, ул. им. Репина, `r`nд. 35/1, кв. 106`r`n"
$ws.Range("B388").Value = "Акционерное общество «Антипинский нефтеперерабатывающий`r`n завод» (АО «Антипинский НПЗ»)"
$ws.Range("C388").Value = "Тюменская область,`r`nг. Тюмень, ул. 6 км Старого Тобольского тракта, д. 20`r`n"
$ws.Range("B389").Value = "Общество с ограниченной ответственностью «Нефтяная компания «Новый Поток» `r`n(ООО «НКНП»)"
$ws.Range("C389").Value = "Оренбургская область, г. Бузулук, ул. Матросова, д. 1`r`n"
$ws.Range("B398").Value = "ООО «Артамира»`r`n"
$ws.Range("C398").Value = "Саратовская обл., г. Саратов, ул. Вольская, д. 91`r`n"
$ws.Range("C400").Value = "Оренбургская область, г. Оренбург,ул. Джангильдина,`r`n д. 3, пом. 12"
$ws.Range("C404").Value = "Самарская область, г. Самара, ул. Ново-Садовая, д. 106, `r`nкорп. 82, оф. 26"
$ws.Range("C405").Value = "Республика Удмуртская, г. Ижевск, ул. им. Репина,`r`nд. 35/1, кв. 106`r`n"
$ws.Range("C406").Value = "Республика Удмуртская, г. Ижевск, ул. им. Репина,`r`nд. 35/1, кв. 106`r`n"
$ws.Range("C407").Value = "Оренбургская обл., г. Оренбург, ул. Комсомольская, `r`nд. 40"
$ws.Range("C408").Value = "Саратовская обл., г. Саратов, ул. им. Пугачева Е.И. , `r`nд. 159, оф. 606`r`n"
$ws.Range("C410").Value = "Ямало-Ненецкий автономный округ, г. Ноябрьск, `r`nул. Ленина,д. 51`r`n"
$ws.Range("C411").Value = "г. Москва, ул. Профсоюзная, д. 56, каб. 77`r`n"
